$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that wraps nothing at the very start
#    of the document (first paragraph "Organización:").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Host "Removed leading _GoBack bookmark"
}

# ---------------------------------------------------------------------------
# 2) Insert two new bullet paragraphs right after the paragraph that ends
#    with "...De esta forma no hay confusiones." (still inside the Trello
#    bulleted list that uses style "Prrafodelista" / numId 1).
# ---------------------------------------------------------------------------
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Las tareas las creamos en grupo al principio de cada sprint*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq 0) {
    throw "Could not locate the 'Las tareas las creamos...' anchor paragraph"
}

# Create a brand-new empty paragraph right after the anchor, then fill that
# *new* paragraph's range via InsertXML (targeting the anchor's own collapsed
# end range would instead overwrite the anchor paragraph's own content).
$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter() | Out-Null

$newParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>miebro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> se añadirá como participante de una tarea al moverla a en proceso.</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Las tareas se podrán dejar inacabadas pero deberán ser pasadas de nuevo a “por hacer”.</w:t></w:r></w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Paragraphs($anchorIndex + 1).Range.InsertXML($newParasXml)
Write-Host "Inserted the two new 'miebro' / 'por hacer' bullet paragraphs"

# ---------------------------------------------------------------------------
# 3) Underline the paragraph mark of the "Crearemos los diagramas..."
#    paragraph (Draw.IO section) and re-anchor the "_GoBack" bookmark at the
#    very end of that paragraph's text.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Crearemos los diagramas*") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the 'Crearemos los diagramas...' paragraph"
}

$targetParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Crearemos los diagramas en los que se basaran las tareas de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trello</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> todos juntos al inicio de cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spring</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($targetParaXml)
Write-Host "Underlined paragraph mark and re-anchored _GoBack bookmark on the Draw.IO paragraph"
